# Stundenliste update: add the "Presentation" row (row 45) to the hours
# log, extending the running-total formula down one row, and update the
# sheet's saved scroll/selection state to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: 2021-09-22, 2 hours, running total, activity "Presentation"
$ws.Range("A45").Value = 44461
$ws.Range("B45").Value = 2
$ws.Range("C45").Formula = "=SUM(C44, B45)"
$ws.Range("D45").Value = "Presentation"

# Restore the view state captured in the saved workbook: scrolled so row
# 19 is at the top, with C44:C45 selected (active cell C44).
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("C44:C45").Select()
